$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the "through" date references
$wb.Worksheets.Item(1).Name = "Through 2021-10-16"
$ws.Range("A11").Value = "October (through 10-16)"

# Update October row (row 11) values for 2016-2021 (columns C-H)
$ws.Range("C11").Value = 29
$ws.Range("D11").Value = 30
$ws.Range("E11").Value = 41
$ws.Range("F11").Value = 20
$ws.Range("G11").Value = 81
$ws.Range("H11").Value = 98

# Update Total row (row 12) values for 2016-2021 (columns C-H)
$ws.Range("C12").Value = 458
$ws.Range("D12").Value = 657
$ws.Range("E12").Value = 589
$ws.Range("F12").Value = 442
$ws.Range("G12").Value = 982
$ws.Range("H12").Value = 1348
